$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: add header F1 = "URL"
$ws.Range("F1").Value = "URL"

# Row 4: new content (was empty A4/B4, linkedin/company F4)
$ws.Range("A4").Value = "Microsoft Jobs, Employment in Denver, CO | Indeed.com"
$ws.Range("B4").Value = "0252655a41544fd28ae41f8b8ff36917@sentry.indeed.com`n"
$ws.Range("C4").Value = 1100
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = 21
$ws.Range("F4").Value = "https://www.indeed.com/q-Microsoft-l-Denver,-CO-jobs.html"

# Row 5: becomes old row6 content w/ updated C/E (old row5 "508,000+..." entry removed)
$ws.Range("A5").Value = "Microsoft Jobs and Careers | Indeed.com"
$ws.Range("B5").Value = "white@2x.png`nwhite@2x.png`n658ba2886a9642c2b8c035add5a02b63@sentry.indeed.com`nu002f9282b91fa44845a98549f9a94b9326b2@sentry.indeed.com`n"
$ws.Range("C5").Value = 282
$ws.Range("D5").Value = 0
$ws.Range("E5").Value = 28
$ws.Range("F5").Value = "https://www.indeed.com/cmp/Microsoft/jobs"

# Row 6: becomes old row4 content (empty A6/B6, linkedin/company F6)
$ws.Range("A6").Value = ""
$ws.Range("B6").Value = ""
$ws.Range("C6").Value = 0
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = 0
$ws.Range("F6").Value = "https://www.linkedin.com/company/microsoft/jobs"

# Row 8: B8 loses the first two lines (donkey@jackasswhisperer.com / u003edonkey@jackasswhisperer.com), counts updated
$ws.Range("B8").Value = "mingraham@theladders.com.`nmingraham@theladders.com`nmlepore@theladders.com.`nmlepore@theladders.com`nudc9e@fashncurious`nmlepore@theladders.com.`nmlepore@theladders.com`neprice@theladders.com`njfabiano@theladders.com`nmlepore@theladders.com.`nmlepore@theladders.com`nu00a0@nytimes`nu00a0@genderfair`nu00a0@janssenglobal`nu00a0@janssenglobal`nmlepore@theladders.com.`nmlepore@theladders.com`nagarrido@theladders.com`nsusannakahr@theladders.com`n"
$ws.Range("C8").Value = 2144
$ws.Range("D8").Value = 35
$ws.Range("E8").Value = 318
